# Appends two more daily "FII DII" blocks (19/06/2024 and 21/06/2024 date
# sections) below the existing data on Sheet1, growing the used range
# from A1:J1022 to A1:J1099. Each block is: a header row (Buying
# Opportunity / support Zone / long buildup / Short buildup / FII
# ENTERING), followed by data rows (ticker symbols in columns A-E and
# numeric stats in columns F-J), closed by a date-stamp row in column A.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the new block (rows 1023-1099, columns A-J) as a 2-D array
# and paste it in a single range write.
$arr = New-Object 'object[,]' 77,10

# row 1023
$arr[0,0] = "Buying Opportunity"
$arr[0,1] = "support Zone"
$arr[0,2] = "long buildup"
$arr[0,3] = "Short buildup"
$arr[0,4] = "FII ENTERING"
# row 1024
$arr[1,0] = "360ONE"
$arr[1,1] = "AAKASH"
$arr[1,2] = "COROMANDEL"
$arr[1,5] = 852.9
$arr[1,6] = 11.93
$arr[1,7] = 1573.1
# row 1025
$arr[2,0] = "AAVAS"
$arr[2,1] = "ARCHIES"
$arr[2,5] = 1935.75
$arr[2,6] = 25.87
# row 1026
$arr[3,0] = "AURIONPRO"
$arr[3,1] = "BHARATRAS"
$arr[3,5] = 2955.65
$arr[3,6] = 11662.35
# row 1027
$arr[4,0] = "AXISBNKETF"
$arr[4,1] = "CCHHL"
$arr[4,5] = 524.02
$arr[4,6] = 17.23
# row 1028
$arr[5,0] = "BANDHANBNK"
$arr[5,1] = "CHEVIOT"
$arr[5,5] = 198.83
$arr[5,6] = 1494.5
# row 1029
$arr[6,0] = "BANKBEES"
$arr[6,1] = "DBCORP"
$arr[6,5] = 526.16
$arr[6,6] = 313.35
# row 1030
$arr[7,0] = "BFSI"
$arr[7,1] = "DIL"
$arr[7,5] = 23.42
$arr[7,6] = 7.75
# row 1031
$arr[8,0] = "GSFC"
$arr[8,1] = "GENSOL"
$arr[8,5] = 239.82
$arr[8,6] = 1015.9
# row 1032
$arr[9,0] = "HDFCBANK"
$arr[9,1] = "GLAXO"
$arr[9,5] = 1657.85
$arr[9,6] = 2563.55
# row 1033
$arr[10,0] = "IDFCFIRSTB"
$arr[10,1] = "GULPOLY"
$arr[10,5] = 82.17
$arr[10,6] = 199.48
# row 1034
$arr[11,0] = "INSECTICID"
$arr[11,1] = "HEROMOTOCO"
$arr[11,5] = 707
$arr[11,6] = 5647.7
# row 1035
$arr[12,0] = "JYOTISTRUC"
$arr[12,1] = "HILTON"
$arr[12,5] = 28.9
$arr[12,6] = 108.43
# row 1036
$arr[13,0] = "LAXMIMACH"
$arr[13,1] = "HNDFDS"
$arr[13,5] = 16620.7
$arr[13,6] = 496.3
# row 1037
$arr[14,0] = "MANGCHEFER"
$arr[14,1] = "IEL"
$arr[14,5] = 123.45
$arr[14,6] = 12.76
# row 1038
$arr[15,0] = "MITCON"
$arr[15,1] = "INDOWIND"
$arr[15,5] = 182.29
$arr[15,6] = 23.73
# row 1039
$arr[16,0] = "MOTILALOFS"
$arr[16,1] = "KAMOPAINTS"
$arr[16,5] = 686.05
$arr[16,6] = 36.9
# row 1040
$arr[17,0] = "NACLIND"
$arr[17,1] = "KIOCL"
$arr[17,5] = 76.1
$arr[17,6] = 452.65
# row 1041
$arr[18,0] = "NFL"
$arr[18,1] = "KIRLPNU"
$arr[18,5] = 125.91
$arr[18,6] = 1285.05
# row 1042
$arr[19,0] = "NPBET"
$arr[19,1] = "KOPRAN"
$arr[19,5] = 261.41
$arr[19,6] = 254.98
# row 1043
$arr[20,0] = "ORIENTELEC"
$arr[20,1] = "KRSNAA"
$arr[20,5] = 248
$arr[20,6] = 642.75
# row 1044
$arr[21,0] = "PSUBNKBEES"
$arr[21,1] = "LINCOLN"
$arr[21,5] = 83.26
$arr[21,6] = 603.4
# row 1045
$arr[22,0] = "PUNJABCHEM"
$arr[22,1] = "LOVABLE"
$arr[22,5] = 1247.4
$arr[22,6] = 124.49
# row 1046
$arr[23,0] = "RCF"
$arr[23,1] = "MAHSCOOTER"
$arr[23,5] = 185.22
$arr[23,6] = 8047.05
# row 1047
$arr[24,0] = "RKDL"
$arr[24,1] = "MUTHOOTFIN"
$arr[24,5] = 28.59
$arr[24,6] = 1740.45
# row 1048
$arr[25,0] = "19/06/2024"
# row 1049
$arr[26,0] = "Buying Opportunity"
$arr[26,1] = "support Zone"
$arr[26,2] = "long buildup"
$arr[26,3] = "Short buildup"
$arr[26,4] = "FII ENTERING"
# row 1050
$arr[27,0] = "AGROPHOS"
$arr[27,1] = "AMBER"
$arr[27,2] = "COROMANDEL"
$arr[27,3] = "MFSL"
$arr[27,4] = "DEEPAKNTR"
$arr[27,5] = 49.91
$arr[27,6] = 3997.45
$arr[27,7] = 1643.8
$arr[27,8] = 987.35
$arr[27,9] = 2605.25
# row 1051
$arr[28,0] = "ASPINWALL"
$arr[28,1] = "APLLTD"
$arr[28,5] = 290.15
$arr[28,6] = 845.3
# row 1052
$arr[29,0] = "ASTEC"
$arr[29,1] = "ARVIND"
$arr[29,5] = 1364.3
$arr[29,6] = 368.1
# row 1053
$arr[30,0] = "BANDHANBNK"
$arr[30,1] = "BEML"
$arr[30,5] = 208.18
$arr[30,6] = 4450.3
# row 1054
$arr[31,0] = "BLUEJET"
$arr[31,1] = "CCHHL"
$arr[31,5] = 424.8
$arr[31,6] = 17.17
# row 1055
$arr[32,0] = "CHAMBLFERT"
$arr[32,1] = "CDSL"
$arr[32,5] = 557.85
$arr[32,6] = 2039.45
# row 1056
$arr[33,0] = "CLEAN"
$arr[33,1] = "DTIL"
$arr[33,5] = 1443.6
$arr[33,6] = 209.2
# row 1057
$arr[34,0] = "COSMOFIRST"
$arr[34,1] = "GLAXO"
$arr[34,5] = 784.85
$arr[34,6] = 2529.85
# row 1058
$arr[35,0] = "DEEPAKFERT"
$arr[35,1] = "IFCI"
$arr[35,5] = 728.15
$arr[35,6] = 62.48
# row 1059
$arr[36,0] = "DEEPAKNTR"
$arr[36,1] = "JCHAC"
$arr[36,5] = 2605.25
$arr[36,6] = 1817.2
# row 1060
$arr[37,0] = "EROSMEDIA"
$arr[37,1] = "JINDRILL"
$arr[37,5] = 20.21
$arr[37,6] = 632.35
# row 1061
$arr[38,0] = "ESSARSHPNG"
$arr[38,1] = "KIRLOSIND"
$arr[38,5] = 58.07
$arr[38,6] = 6081.7
# row 1062
$arr[39,0] = "GSFC"
$arr[39,1] = "KPIGREEN"
$arr[39,5] = 262.66
$arr[39,6] = 1778.9
# row 1063
$arr[40,0] = "INOXGREEN"
$arr[40,1] = "LICI"
$arr[40,5] = 159.32
$arr[40,6] = 1027.4
# row 1064
$arr[41,0] = "JYOTISTRUC"
$arr[41,1] = "METROPOLIS"
$arr[41,5] = 29.77
$arr[41,6] = 1946.35
# row 1065
$arr[42,0] = "KERNEX"
$arr[42,1] = "PTC"
$arr[42,5] = 447.7
$arr[42,6] = 206.32
# row 1066
$arr[43,0] = "KPRMILL"
$arr[43,1] = "PTCIL"
$arr[43,5] = 890.75
$arr[43,6] = 13603.15
# row 1067
$arr[44,0] = "MADRASFERT"
$arr[44,1] = "SAHYADRI"
$arr[44,5] = 128.73
$arr[44,6] = 415.45
# row 1068
$arr[45,0] = "MANGCHEFER"
$arr[45,5] = 134.54
# row 1069
$arr[46,0] = "NAVINFLUOR"
$arr[46,5] = 3770.25
# row 1070
$arr[47,0] = "NFL"
$arr[47,5] = 150.91
# row 1071
$arr[48,0] = "OPTIEMUS"
$arr[48,5] = 310.45
# row 1072
$arr[49,0] = "PARADEEP"
$arr[49,5] = 92.37
# row 1073
$arr[50,0] = "PKTEA"
$arr[50,5] = 344.5
# row 1074
$arr[51,0] = "PUNJABCHEM"
$arr[51,5] = 1370
# row 1075
$arr[52,0] = "RAMKY"
$arr[52,5] = 594.9
# row 1076
$arr[53,0] = "RCF"
$arr[53,5] = 222.26
# row 1077
$arr[54,0] = "RKDL"
$arr[54,5] = 30.01
# row 1078
$arr[55,0] = "RTNINDIA"
$arr[55,5] = 81.83
# row 1079
$arr[56,0] = "21/06/2024"
# row 1080
$arr[57,0] = "Buying Opportunity"
$arr[57,1] = "support Zone"
$arr[57,2] = "long buildup"
$arr[57,3] = "Short buildup"
$arr[57,4] = "FII ENTERING"
# row 1081
$arr[58,0] = "ADFFOODS"
$arr[58,1] = "AETHER"
$arr[58,5] = 243.38
$arr[58,6] = 863.2
# row 1082
$arr[59,0] = "FCL"
$arr[59,1] = "AGARIND"
$arr[59,5] = 390.15
$arr[59,6] = 1151.6
# row 1083
$arr[60,0] = "GMBREW"
$arr[60,1] = "ALPHAETF"
$arr[60,5] = 651.35
$arr[60,6] = 27.41
# row 1084
$arr[61,0] = "IKIO"
$arr[61,1] = "AMNPLST"
$arr[61,5] = 316.4
$arr[61,6] = 251.51
# row 1085
$arr[62,0] = "ITETF"
$arr[62,1] = "ANUP"
$arr[62,5] = 35.77
$arr[62,6] = 1910.05
# row 1086
$arr[63,0] = "IVZINGOLD"
$arr[63,1] = "CONCOR"
$arr[63,5] = 6463.55
$arr[63,6] = 1090.9
# row 1087
$arr[64,0] = "KANORICHEM"
$arr[64,1] = "CPSEETF"
$arr[64,5] = 130.78
$arr[64,6] = 93.18
# row 1088
$arr[65,0] = "KPITTECH"
$arr[65,1] = "DYNAMATECH"
$arr[65,5] = 1588.75
$arr[65,6] = 7727.55
# row 1089
$arr[66,0] = "LGBFORGE"
$arr[66,1] = "GANESHHOUC"
$arr[66,5] = 14.3
$arr[66,6] = 959.6
# row 1090
$arr[67,0] = "LPDC"
$arr[67,1] = "HARDWYN"
$arr[67,5] = 9.87
$arr[67,6] = 31.84
# row 1091
$arr[68,0] = "MIRZAINT"
$arr[68,1] = "HINDMOTORS"
$arr[68,5] = 46.02
$arr[68,6] = 35.16
# row 1092
$arr[69,0] = "OMAXE"
$arr[69,1] = "IFCI"
$arr[69,5] = 109.3
$arr[69,6] = 61.18
# row 1093
$arr[70,0] = "PCBL"
$arr[70,1] = "MMTC"
$arr[70,5] = 266.2
$arr[70,6] = 80.94
# row 1094
$arr[71,0] = "PKTEA"
$arr[71,1] = "MOMENTUM"
$arr[71,5] = 413.4
$arr[71,6] = 35.81
# row 1095
$arr[72,0] = "RHL"
$arr[72,1] = "NDRAUTO"
$arr[72,5] = 206.7
$arr[72,6] = 793.8
# row 1096
$arr[73,0] = "RSWM"
$arr[73,1] = "NECLIFE"
$arr[73,5] = 198.71
$arr[73,6] = 32.92
# row 1097
$arr[74,1] = "ORISSAMINE"
$arr[74,6] = 7451.1
# row 1098
$arr[75,1] = "RUCHINFRA"
$arr[75,6] = 12.34
# row 1099
$arr[76,0] = "21/06/2024"

$ws.Range("A1023:J1099").Value = $arr
